# Updated version of document that includes summer 2020 field site data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert 16 new rows at row 17 (pushes existing rows 17-38 down to 33-54),
# inheriting the formatting of the row above (style index 1).
$ws.Range("A17:A32").EntireRow.Insert()

# Row 17 - first new site; plain (non-shared) formula
$ws.Range("A17").Value = "DAYSCH01"
$ws.Range("B17").Value = "S Umpqua"
$ws.Range("C17").Value = 129912.24
$ws.Range("D17").Formula = "=C17/1000"

# Rows 18-32 - remaining new sites; fill column A/B/C per-row, then set the
# D-column formula across the whole block in one go so Excel compresses it
# into a single shared formula (matching the authored diff's t="shared").
$ws.Range("A18").Value = "DAYSCH02"
$ws.Range("B18").Value = "S Umpqua"
$ws.Range("C18").Value = 129960.838

$ws.Range("A19").Value = "DAYSCH03"
$ws.Range("B19").Value = "S Umpqua"
$ws.Range("C19").Value = 129942.77099999999

$ws.Range("A20").Value = "DAYSCH04"
$ws.Range("B20").Value = "S Umpqua"
$ws.Range("C20").Value = 129995.481

$ws.Range("A21").Value = "DAYSCH05"
$ws.Range("B21").Value = "S Umpqua"
$ws.Range("C21").Value = 130182.427

$ws.Range("A22").Value = "DAYSCH06"
$ws.Range("B22").Value = "S Umpqua"
$ws.Range("C22").Value = 130316.891

$ws.Range("A23").Value = "DAYSCH07"
$ws.Range("B23").Value = "S Umpqua"
$ws.Range("C23").Value = 130128.48699999999

$ws.Range("A24").Value = "DAYSCH08"
$ws.Range("B24").Value = "S Umpqua"
$ws.Range("C24").Value = 130094.93

$ws.Range("A25").Value = "SR01"
$ws.Range("B25").Value = "S Umpqua"
$ws.Range("C25").Value = 157348.109

$ws.Range("A26").Value = "SR02"
$ws.Range("B26").Value = "S Umpqua"
$ws.Range("C26").Value = 157580.96599999999

$ws.Range("A27").Value = "OB01"
$ws.Range("B27").Value = "S Umpqua"
$ws.Range("C27").Value = 159871.80499999999

$ws.Range("A28").Value = "OB02"
$ws.Range("B28").Value = "S Umpqua"
$ws.Range("C28").Value = 159945.03599999999

$ws.Range("A29").Value = "OB03"
$ws.Range("B29").Value = "S Umpqua"
$ws.Range("C29").Value = 159984.503

$ws.Range("A30").Value = "OB04"
$ws.Range("B30").Value = "S Umpqua"
$ws.Range("C30").Value = 160109.67300000001

$ws.Range("A31").Value = "ZINCCMP01"
$ws.Range("B31").Value = "S Umpqua"
$ws.Range("C31").Value = 196288.24299999999

$ws.Range("A32").Value = "ZINCCMP02"
$ws.Range("B32").Value = "S Umpqua"
$ws.Range("C32").Value = 196196.88800000001

$ws.Range("D18:D32").Formula = "=C18/1000"

# The old final row (previously row 38, "ZERO_UmpConf") is no longer part
# of this dataset; after the insert above it now lives at row 54 - remove it.
$ws.Range("A54").EntireRow.Delete()

# Leave the cell cursor where the author left it when they saved.
[void]$ws.Range("C33").Select()
